# Apply the "completed requirements matrix" edit to the Checklist sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# ---------------------------------------------------------------------
# 1. Un-hide the detail rows for sections 1 and 2 (rows 3-25), and clear
#    the "collapsed" state of their group header rows (2 and 8 / 26).
# ---------------------------------------------------------------------
for ($r = 3; $r -le 25; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# ---------------------------------------------------------------------
# 2. Update the "Database files" answer (requirement 3.4, cell C30):
#    seed-data.sql -> users.sql, and the load-schema.sh list entry ->
#    load-users.sh (only the bullet list item, not the prose below it
#    that still references load-schema.sh).
# ---------------------------------------------------------------------
$c30 = $ws.Range("C30")
$newText = $c30.Value2.Replace("    - seed-data.sql`n", "    - users.sql`n")
$newText = $newText.Replace("    - load-schema.sh`n", "    - load-users.sh`n")
$c30.Value = $newText

# ---------------------------------------------------------------------
# 3. Mark section "3" (Other Requirements and Guidelines) as complete:
#    style A26 like a completed row and add a "complete" note in C26.
# ---------------------------------------------------------------------
$ws.Range("A27").Copy()
$ws.Range("A26").PasteSpecial(-4122)

$ws.Range("C26").Value = "Other requirements 3.1-3-7 complete"
$ws.Range("C29").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Rows.Item(26).RowHeight = 22

# ---------------------------------------------------------------------
# 4. Mark requirements 3.6.1 - 3.6.5 (rows 33-37) as complete: style the
#    A-column like other completed rows and add "check" answers in C.
# ---------------------------------------------------------------------
$ws.Range("A27").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A37").PasteSpecial(-4122)

$ws.Range("C33").Value = "check"
$ws.Range("C34").Value = "check"
$ws.Range("C35").Value = "check"
$ws.Range("C36").Value = "check"
$ws.Range("C37").Value = "check"

$ws.Range("C29").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C35").PasteSpecial(-4122)
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C37").PasteSpecial(-4122)

$ws.Rows.Item(35).RowHeight = 22
$ws.Rows.Item(36).RowHeight = 22
$ws.Rows.Item(37).RowHeight = 22

# ---------------------------------------------------------------------
# 5. Update the view state: scroll position and active selection.
# ---------------------------------------------------------------------
$ws.Range("C27").Select()
